# Add the "Change Column Type" worksheet (GOMS task) after
# "Principal Component Analysis", populate it with the same
# Action / Time / Content layout used by the other task sheets,
# and make it the active sheet/tab.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Change Column Type"

# Header row
$ws.Range("A1").Value = "Action"
$ws.Range("B1").Value = "Time"
$ws.Range("C1").Value = "Content"

# Upload CSV
$ws.Range("A2").Value = "Upload CSV"
$ws.Range("B2").Value = "5 min"
$ws.Range("C2").Value = "df = pd.read_csv('file.csv')"

# Identify Column & Type
$ws.Range("A3").Value = "Identify Column & Type"
$ws.Range("B3").Value = "2 min"
$ws.Range("C3").Value = "Review the data with df.dtypes"

# Convert Data Type
$ws.Range("A4").Value = "Convert Data Type"
$ws.Range("B4").Value = "2 min"
$ws.Range("C4").Value = "df['column'] = df['column'].astype('desired_type')"

# Verify Changes
$ws.Range("A5").Value = "Verify Changes"
$ws.Range("B5").Value = "1 min"
$ws.Range("C5").Value = "df.dtypes to confirm the change"

# Overall
$ws.Range("A6").Value = "Overall"
$ws.Range("B6").Value = "10 min"

# Formatting to match the rest of the workbook's task sheets:
# 13pt body text, bold 13pt header row and bold "Overall" row.
$ws.Range("A1:C5").Font.Size = 13
$ws.Range("A6:B6").Font.Size = 13
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A6:B6").Font.Bold = $true

$ws.Rows.Item(1).RowHeight = 17
$ws.Rows.Item(2).RowHeight = 17
$ws.Rows.Item(3).RowHeight = 17
$ws.Rows.Item(4).RowHeight = 17
$ws.Rows.Item(5).RowHeight = 17
$ws.Rows.Item(6).RowHeight = 17

# Matches the saved selection/active cell on the new tab.
$ws.Range("O18").Select() | Out-Null
